## "add text and integer inputs with do not know button"
##
## Adds two new rows to the "survey" sheet describing two new custom
## prompt types (custom_number "number_3" and custom_text "text_1"),
## each with a new inputAttributes.showIDK = TRUE flag, registers the
## new "custom_text" prompt type on the "prompt_types" sheet, and
## repoints the existing "text" question (row 3) at the new
## "custom_text" prompt type.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("survey")
$ws3 = $wb.Worksheets.Item("prompt_types")

## --- survey sheet -------------------------------------------------

# Make room for the two new rows right after the existing
# "custom_number" example (old row 11 "end screen" shifts to row 13).
$ws1.Rows("11:12").Insert() | Out-Null

# Fill in the new cells roughly in the order a person filling out the
# form left-to-right / row-by-row would hit them.
$ws1.Range("E11").Value() = "number_3"
$ws1.Range("H11").Value() = "Same as above, but with additional button for don't know"

$ws1.Range("K1").Value()  = "inputAttributes.showIDK"

# The plain "text" prompt becomes the new "custom_text" prompt type.
$ws1.Range("C3").Value() = "custom_text"

$ws1.Range("E12").Value() = "text_1"
$ws1.Range("F12").Value() = "Input text, or don't know (-99)"
$ws1.Range("H12").Value() = "Same as above for text input"

$ws1.Range("F11").Value() = "Input a number, or don't know (-99)"

$ws1.Range("C11").Value() = "custom_number"
$ws1.Range("C12").Value() = "custom_text"

$ws1.Range("K11").Value() = $true
$ws1.Range("K12").Value() = $true

# Row heights for the two new rows (content wraps over multiple lines).
$ws1.Rows.Item(11).RowHeight = 38.25
$ws1.Rows.Item(12).RowHeight = 25.5

## --- prompt_types sheet --------------------------------------------

# Register the new custom_text prompt type (maps to the base "string"
# xlsform type, same pattern as the existing custom_number/integer and
# custom_date/text rows above it).
$ws3.Range("A5").Value() = "string"
$ws3.Range("B5").Value() = "custom_text"

# Bold the header row to call out the prompt type columns.
$ws3.Range("A1:B1").Font.Bold = $true

$ws3.PageSetup.PaperSize = 9
$ws3.PageSetup.Orientation = 1

## --- view / selection bookkeeping -----------------------------------

$ws3.Range("B11").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("C3").Select() | Out-Null
$excel.ActiveWindow.Zoom = 115
